$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data refresh adds 4 new price observations (for the new
# reporting date 2023-11-09, serial 45239) at the top of the existing
# date-ordered block for this market/variety, pushing the previously
# existing rows (old 1205..1230) down to 1209..1234.
$ws.Rows("1205:1208").Insert()

# Row 1205: Repollo - Copenhague
$ws.Cells.Item(1205, 1).Value = 10
$ws.Cells.Item(1205, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(1205, 3).Value = "La Araucanía"
$ws.Cells.Item(1205, 4).Value = 45239
$ws.Cells.Item(1205, 5).Value = 9
$ws.Cells.Item(1205, 6).Value = 100112006
$ws.Cells.Item(1205, 7).Value = "Repollo"
$ws.Cells.Item(1205, 8).Value = "Copenhague"
$ws.Cells.Item(1205, 9).Value = "Primera"
$ws.Cells.Item(1205, 10).Value = 1500
$ws.Cells.Item(1205, 11).Value = 1400
$ws.Cells.Item(1205, 12).Value = 1400
$ws.Cells.Item(1205, 13).Value = 1400
$ws.Cells.Item(1205, 14).Value = "`$/unidad"
$ws.Cells.Item(1205, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(1205, 16).Value = 1400
$ws.Cells.Item(1205, 17).Value = 1
$ws.Cells.Item(1205, 18).Value = "Hortaliza"

# Row 1206: Repollo - Crespo record
$ws.Cells.Item(1206, 1).Value = 10
$ws.Cells.Item(1206, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(1206, 3).Value = "La Araucanía"
$ws.Cells.Item(1206, 4).Value = 45239
$ws.Cells.Item(1206, 5).Value = 9
$ws.Cells.Item(1206, 6).Value = 100112006
$ws.Cells.Item(1206, 7).Value = "Repollo"
$ws.Cells.Item(1206, 8).Value = "Crespo record"
$ws.Cells.Item(1206, 9).Value = "Primera"
$ws.Cells.Item(1206, 10).Value = 1000
$ws.Cells.Item(1206, 11).Value = 1400
$ws.Cells.Item(1206, 12).Value = 1500
$ws.Cells.Item(1206, 13).Value = 1440
$ws.Cells.Item(1206, 14).Value = "`$/unidad"
$ws.Cells.Item(1206, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1206, 16).Value = 1440
$ws.Cells.Item(1206, 17).Value = 1
$ws.Cells.Item(1206, 18).Value = "Hortaliza"

# Row 1207: Repollo - Crespo record
$ws.Cells.Item(1207, 1).Value = 10
$ws.Cells.Item(1207, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(1207, 3).Value = "La Araucanía"
$ws.Cells.Item(1207, 4).Value = 45239
$ws.Cells.Item(1207, 5).Value = 9
$ws.Cells.Item(1207, 6).Value = 100112006
$ws.Cells.Item(1207, 7).Value = "Repollo"
$ws.Cells.Item(1207, 8).Value = "Crespo record"
$ws.Cells.Item(1207, 9).Value = "Primera"
$ws.Cells.Item(1207, 10).Value = 3300
$ws.Cells.Item(1207, 11).Value = 1400
$ws.Cells.Item(1207, 12).Value = 1500
$ws.Cells.Item(1207, 13).Value = 1455
$ws.Cells.Item(1207, 14).Value = "`$/unidad"
$ws.Cells.Item(1207, 15).Value = "Región del Maule"
$ws.Cells.Item(1207, 16).Value = 1455
$ws.Cells.Item(1207, 17).Value = 1
$ws.Cells.Item(1207, 18).Value = "Hortaliza"

# Row 1208: Repollo - Morada(o)
$ws.Cells.Item(1208, 1).Value = 10
$ws.Cells.Item(1208, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(1208, 3).Value = "La Araucanía"
$ws.Cells.Item(1208, 4).Value = 45239
$ws.Cells.Item(1208, 5).Value = 9
$ws.Cells.Item(1208, 6).Value = 100112006
$ws.Cells.Item(1208, 7).Value = "Repollo"
$ws.Cells.Item(1208, 8).Value = "Morada(o)"
$ws.Cells.Item(1208, 9).Value = "Primera"
$ws.Cells.Item(1208, 10).Value = 500
$ws.Cells.Item(1208, 11).Value = 1500
$ws.Cells.Item(1208, 12).Value = 1500
$ws.Cells.Item(1208, 13).Value = 1500
$ws.Cells.Item(1208, 14).Value = "`$/unidad"
$ws.Cells.Item(1208, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(1208, 16).Value = 1500
$ws.Cells.Item(1208, 17).Value = 1
$ws.Cells.Item(1208, 18).Value = "Hortaliza"
